$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so values like "1.013" or
# "20.588.59" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.588.59"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.477.58"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "0.9580"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").Value = "280.37"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "0.3661"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("D8").Value = "0.3069"
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "40.06"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "1.062"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "0.06679"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "5.527"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "18.07"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "0.9592"
$ws.Range("E16").Value = "  +5.71%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "1.478.75"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "0.05951"
$ws.Range("E19").Value = "  +4.10%  "
$ws.Range("D20").Value = "70.00"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "5.506"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").Value = "14.45"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "2.260"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "20.633.52"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "143.26"
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("D27").Value = "2.111"
$ws.Range("E27").Value = "  -8.25%  "
$ws.Range("D28").Value = "17.29"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "1.639.27"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "113.92"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "3.975"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "5.032"
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").Value = "0.8136"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").Value = "0.07970"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "1.532"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "1.214"
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").Value = "0.05812"
$ws.Range("E37").Value = "  -5.39%  "
$ws.Range("D38").Value = "4.758"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "0.02048"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "0.9602"
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("D41").Value = "10.39"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "7.458"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").Value = "0.5314"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "12.30"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "118.31"
$ws.Range("E47").Value = "  -4.90%  "
$ws.Range("D48").Value = "0.5204"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "1.825"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "0.06493"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "0.9891"
$ws.Range("E51").Value = "  -0.32%  "
